# Auto-generated: apply scheduled-runner market data refresh to Belias_Profits workbook.
# Updates cached price/profit columns (H:N) on specific rows across multiple sheets.
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 2850.5256
$ws.Range("I15").Value = 2850.5256
$ws.Range("K15").Value = 8551.576799999999
$ws.Range("M15").Value = -8382.576799999999
$ws.Range("H18").Value = 829.381
$ws.Range("I18").Value = 794.875
$ws.Range("J18").Value = 939.8
$ws.Range("K18").Value = 794.875
$ws.Range("L18").Value = 939.8
$ws.Range("M18").Value = -510.875
$ws.Range("N18").Value = -1507.8
$ws.Range("H98").Value = 2555.8823
$ws.Range("I98").Value = 1541.25
$ws.Range("J98").Value = 4991
$ws.Range("K98").Value = 1541.25
$ws.Range("L98").Value = 4991
$ws.Range("M98").Value = -43.25
$ws.Range("N98").Value = -7987
$ws.Range("H116").Value = 2101.9092
$ws.Range("I116").Value = 2002.8572
$ws.Range("J116").Value = 2275.25
$ws.Range("K116").Value = 2002.8572
$ws.Range("L116").Value = 2275.25
$ws.Range("M116").Value = 1439.1428
$ws.Range("N116").Value = -9159.25
$ws.Range("H122").Value = 2555.8823
$ws.Range("I122").Value = 1541.25
$ws.Range("J122").Value = 4991
$ws.Range("K122").Value = 4623.75
$ws.Range("L122").Value = 14973
$ws.Range("M122").Value = -2173.75
$ws.Range("N122").Value = -19873
$ws.Range("H135").Value = 47619364
$ws.Range("I135").Value = 292.33334
$ws.Range("J135").Value = 333333800
$ws.Range("K135").Value = 2631.00006
$ws.Range("L135").Value = 3000004200
$ws.Range("M135").Value = -96.0000600000003
$ws.Range("N135").Value = -3000009270
$ws.Range("H138").Value = 20591.166
$ws.Range("I138").Value = 1454
$ws.Range("J138").Value = 55843.844
$ws.Range("K138").Value = 4362
$ws.Range("L138").Value = 167531.532
$ws.Range("M138").Value = 778
$ws.Range("N138").Value = -177811.532

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 6995827.5
$ws.Range("I45").Value = 10103529
$ws.Range("J45").Value = 3500
$ws.Range("K45").Value = 10103529
$ws.Range("L45").Value = 3500
$ws.Range("M45").Value = -10103152
$ws.Range("N45").Value = -4254

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 641.9167
$ws.Range("I80").Value = 998.41174
$ws.Range("J80").Value = 322.94736
$ws.Range("K80").Value = 998.41174
$ws.Range("L80").Value = 322.94736
$ws.Range("M80").Value = -0.4117400000000089
$ws.Range("N80").Value = -2318.94736
$ws.Range("H83").Value = 641.9167
$ws.Range("I83").Value = 998.41174
$ws.Range("J83").Value = 322.94736
$ws.Range("K83").Value = 4992.0587
$ws.Range("L83").Value = 1614.7368
$ws.Range("M83").Value = -0.05869999999958964
$ws.Range("N83").Value = -11598.7368
$ws.Range("H105").Value = 1877.8695
$ws.Range("I105").Value = 1705.625
$ws.Range("J105").Value = 2271.5715
$ws.Range("K105").Value = 1705.625
$ws.Range("L105").Value = 2271.5715
$ws.Range("M105").Value = 41.375
$ws.Range("N105").Value = -5765.5715
$ws.Range("H108").Value = 25500
$ws.Range("J108").Value = 25500
$ws.Range("L108").Value = 25500
$ws.Range("N108").Value = -33180
$ws.Range("H133").Value = 36612.168
$ws.Range("J133").Value = 36612.168
$ws.Range("L133").Value = 36612.168
$ws.Range("N133").Value = -46732.168

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H106").Value = 25000
$ws.Range("J106").Value = 25000
$ws.Range("L106").Value = 25000
$ws.Range("N106").Value = -27524
$ws.Range("H134").Value = 1735.7949
$ws.Range("I134").Value = 1224.875
$ws.Range("J134").Value = 2553.2666
$ws.Range("K134").Value = 3674.625
$ws.Range("L134").Value = 7659.7998
$ws.Range("M134").Value = -1139.625
$ws.Range("N134").Value = -12729.7998

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 812.81396
$ws.Range("I4").Value = 111.61539
$ws.Range("J4").Value = 1116.6666
$ws.Range("K4").Value = 334.84617
$ws.Range("L4").Value = 3349.9998
$ws.Range("M4").Value = -222.84617
$ws.Range("N4").Value = -3573.9998
$ws.Range("H7").Value = 50.625
$ws.Range("I7").Value = 49.11111
$ws.Range("J7").Value = 52.57143
$ws.Range("K7").Value = 147.33333
$ws.Range("L7").Value = 157.71429
$ws.Range("M7").Value = -35.33332999999999
$ws.Range("N7").Value = -381.71429

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1222.3914
$ws.Range("I7").Value = 1214.619
$ws.Range("J7").Value = 1304
$ws.Range("K7").Value = 1214.619
$ws.Range("L7").Value = 1304
$ws.Range("M7").Value = -1102.619
$ws.Range("N7").Value = -1528
$ws.Range("H40").Value = 1435.7446
$ws.Range("I40").Value = 1483.0952
$ws.Range("J40").Value = 1038
$ws.Range("K40").Value = 1483.0952
$ws.Range("L40").Value = 1038
$ws.Range("M40").Value = -1347.0952
$ws.Range("N40").Value = -1310
$ws.Range("H82").Value = 1055.1666
$ws.Range("I82").Value = 902
$ws.Range("J82").Value = 1074.3125
$ws.Range("K82").Value = 902
$ws.Range("L82").Value = 1074.3125
$ws.Range("M82").Value = -541
$ws.Range("N82").Value = -1796.3125
$ws.Range("H85").Value = 1055.1666
$ws.Range("I85").Value = 902
$ws.Range("J85").Value = 1074.3125
$ws.Range("K85").Value = 902
$ws.Range("L85").Value = 1074.3125
$ws.Range("M85").Value = 346
$ws.Range("N85").Value = -3570.3125
$ws.Range("H104").Value = 20720.941
$ws.Range("J104").Value = 20720.941
$ws.Range("L104").Value = 20720.941
$ws.Range("N104").Value = -27708.941
$ws.Range("H126").Value = 1222.3914
$ws.Range("I126").Value = 1214.619
$ws.Range("J126").Value = 1304
$ws.Range("K126").Value = 3643.857
$ws.Range("L126").Value = 3912
$ws.Range("M126").Value = -1173.857
$ws.Range("N126").Value = -8852

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 151663.12
$ws.Range("I122").Value = 1266.6666
$ws.Range("J122").Value = 241901
$ws.Range("K122").Value = 3799.9998
$ws.Range("L122").Value = 725703
$ws.Range("M122").Value = -1349.9998
$ws.Range("N122").Value = -730603
$ws.Range("H126").Value = 1600.3077
$ws.Range("I126").Value = 1464
$ws.Range("J126").Value = 2350
$ws.Range("K126").Value = 4392
$ws.Range("L126").Value = 7050
$ws.Range("M126").Value = -1922
$ws.Range("N126").Value = -11990
